$d = $word.ActiveDocument

# --- Paragraph 1 ("Midterm exam" title): bump font size to 20pt (w:sz/w:szCs = 40) ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Font.Size = 20
$r1.Font.SizeBi = 20

# --- Move the _GoBack bookmark from paragraph 2 to paragraph 3, and give ---
# --- paragraph 3 a 12pt (w:sz/w:szCs = 24) paragraph-mark run size.      ---

# Give the (currently empty) 3rd paragraph a run temporarily so the size change
# actually sticks to its paragraph-mark run properties, then delete the inserted
# text (but not the paragraph mark) so the paragraph stays empty.
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$insertStart = $r3.Start
$r3.InsertAfter("X")
$r3.Font.Size = 12
$r3.Font.SizeBi = 12
$tmp = $d.Range($insertStart, $insertStart + 1)
$tmp.Delete()

# Remove the bookmark from paragraph 2 and re-add it at paragraph 3.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$p3 = $d.Paragraphs.Item(3)
$d.Bookmarks.Add("_GoBack", $p3.Range)
